$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typos in recipe directions / names
$ws.Range("C2").Value = "Put all in a rice cooker or cook quinoa and veggies separately then add everything together"
$ws.Range("C6").Value = "Pour broth in pot, cut up veggies and put in broth, put chicken in broth, bring to boil and put noodles in, season to taste, boil until chicken and noodles cooked, pull apart chicken and put back in soup"
$ws.Range("C7").Value = "Create egg and cream mixture with spices, spray pan, thin walled baking pan preferred bake at 325 for 25-30 minutes"
$ws.Range("A8").Value = "Zucchini Bread"
$ws.Range("C9").Value = "Oven to 400F, baking sheet with parchment paper with oil on top, combine ingredients and let rest for 15 mins, bake for 30 mins, flip halfway thru"
$ws.Range("C10").Value = "Pre cook onions, wash and dry beans, then mix everything in pan and cook with cast iron skillet if possible"
$ws.Range("A12").Value = "Tempeh Ratatouille"
$ws.Range("C15").Value = "Get large soup pan, cut up veggies and simmer in oil till soft, add broth and chicken simmer for 30 mins, pull apart chicken with forks, return chicken to pot, stir in parsley and spinach and lemon juice"

# Update the selected cell in the sheet view
$ws.Range("C23").Select()

# Adjust row 15 height slightly
$ws.Rows.Item(15).RowHeight = 23.85
